# Update the "想去人数" (F column) counts on the sheets that contain the
# event data ("展览" and "全部类型"). The "演出" and "本地生活" sheets only
# contain a header row and are left untouched.

$wb = $excel.ActiveWorkbook

# Row number (as in the sheet) -> new value for column F
$updates = @{
    3  = 99
    4  = 280
    6  = 579
    7  = 61
    8  = 2049
    9  = 69
    11 = 4449
    13 = 284
    14 = 105
    15 = 11
    16 = 123
    19 = 78
    20 = 3289
    22 = 496
    24 = 17
    25 = 80
    26 = 92
    27 = 11
    29 = 59
    30 = 202
    32 = 618
    33 = 1899
    34 = 295
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
